# Refresh the NATMI ligand-receptor TPM output (Dkk2-Kremen2, YoungD7) with the
# re-run numbers: a new "Resolving-Mac" target cluster is now included, so every
# sending cluster (ECs, FAPs, MuSCs) gets a 3rd row, and the expression-derived
# metrics for the existing ECs/MuSCs target rows shift accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol | D Target cluster
#          E..J  Ligand-side stats | K..P Receptor-side stats | Q..T Edge-weight stats
$data = @(
        @("ECs","Dkk2","Kremen2","ECs",2,0.6666666666666666,0.2928123333333333,0.8784369999999999,0.06406943071632207,0.06406943071632207,2,0.6666666666666666,0.3987243333333333,1.196173,0.5027922793035905,0.5027922793035905,0.1167514024001111,1.050762621601,0.03221361510354304,0.03221361510354304),
        @("ECs","Dkk2","Kremen2","MuSCs",2,0.6666666666666666,0.2928123333333333,0.8784369999999999,0.06406943071632207,0.06406943071632207,2,0.6666666666666666,0.321599,0.964797,0.4055370608559684,0.4055370608559684,0.09416815358766666,0.847513382289,0.02598252862341235,0.02598252862341235),
        @("ECs","Dkk2","Kremen2","Resolving-Mac",2,0.6666666666666666,0.2928123333333333,0.8784369999999999,0.06406943071632207,0.06406943071632207,1,0.3333333333333333,0.07269666666666667,0.21809,0.09167065984044119,0.09167065984044119,0.02128648059222222,0.19157832533,0.005873286989366675,0.005873286989366675),
        @("FAPs","Dkk2","Kremen2","ECs",3,1,3.944536333333333,11.833609,0.8630927339690215,0.8630927339690215,2,0.6666666666666666,0.3987243333333333,1.196173,0.5027922793035905,0.5027922793035905,1.572782619817444,14.155043578357,0.4339563629626518,0.4339563629626518),
        @("FAPs","Dkk2","Kremen2","MuSCs",3,1,3.944536333333333,11.833609,0.8630927339690215,0.8630927339690215,2,0.6666666666666666,0.321599,0.964797,0.4055370608559684,0.4055370608559684,1.268558940263667,11.417030462373,0.3500160905799392,0.3500160905799392),
        @("FAPs","Dkk2","Kremen2","Resolving-Mac",3,1,3.944536333333333,11.833609,0.8630927339690215,0.8630927339690215,1,0.3333333333333333,0.07269666666666667,0.21809,0.09167065984044119,0.09167065984044119,0.2867546429788889,2.58079178681,0.07912028042643057,0.07912028042643057),
        @("MuSCs","Dkk2","Kremen2","ECs",3,1,0.332886,0.998658,0.07283783531465635,0.07283783531465635,2,0.6666666666666666,0.3987243333333333,1.196173,0.5027922793035905,0.5027922793035905,0.132729748426,1.194567735834,0.03662230123739562,0.03662230123739562),
        @("MuSCs","Dkk2","Kremen2","MuSCs",3,1,0.332886,0.998658,0.07283783531465635,0.07283783531465635,2,0.6666666666666666,0.321599,0.964797,0.4055370608559684,0.4055370608559684,0.107055804714,0.963502242426,0.02953844165261679,0.02953844165261679),
        @("MuSCs","Dkk2","Kremen2","Resolving-Mac",3,1,0.332886,0.998658,0.07283783531465635,0.07283783531465635,1,0.3333333333333333,0.07269666666666667,0.21809,0.09167065984044119,0.09167065984044119,0.02419970258,0.21779732322,0.006677092424643936,0.006677092424643936)
    )

# Write row-by-row / cell-by-cell starting at A2 (row 1 is the header).
# This also extends the used range to A1:T10 and adds "Resolving-Mac" to the
# shared-string table the first time it's assigned to a cell.
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($startRow + $i, $j + 1).Value = $row[$j]
    }
}
